$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# ---------------------------------------------------------------------------
# Shift the "sum" row (19) and everything below it down by one row. This
# turns row 19 -> 20 (sum row), row 25 -> 26 (initial balance) and
# row 27 -> 28 (balance =), leaving row 19 blank just like in the target.
# ---------------------------------------------------------------------------
$ws.Rows.Item(19).Insert()

# ---------------------------------------------------------------------------
# Row 13: add the new "considered" columns L and P.
# ---------------------------------------------------------------------------
$ws.Range("L13").Formula = '=IF(K13<>"",K13,J13)'
$ws.Range("P13").Formula = '=IF(O13<>"",O13,N13)'

# ---------------------------------------------------------------------------
# Row 14: add L14/P14; K14/O14 keep their existing formulas.
# ---------------------------------------------------------------------------
$ws.Range("L14").Formula = '=IF(K14<>"",K14,J14)'
$ws.Range("P14").Formula = '=IF(O14<>"",O14,N14)'

# ---------------------------------------------------------------------------
# Row 15: add L15/P15.
# ---------------------------------------------------------------------------
$ws.Range("L15").Formula = '=IF(K15<>"",K15,J15)'
$ws.Range("P15").Formula = '=IF(O15<>"",O15,N15)'

# ---------------------------------------------------------------------------
# Row 16: J16 is removed, K16 becomes a hard-coded value (no longer part of
# the shared formula chain), L16/P16 are added, N16/O16 stay as-is.
# ---------------------------------------------------------------------------
$ws.Range("J16").ClearContents()
$ws.Range("K16").Value = 10.130000000000001
$ws.Range("L16").Formula = '=IF(K16<>"",K16,J16)'
$ws.Range("P16").Formula = '=IF(O16<>"",O16,N16)'

# ---------------------------------------------------------------------------
# Row 17: add L17/P17.
# ---------------------------------------------------------------------------
$ws.Range("L17").Formula = '=IF(K17<>"",K17,J17)'
$ws.Range("P17").Formula = '=IF(O17<>"",O17,N17)'

# ---------------------------------------------------------------------------
# Row 18: K18, N18 and O18 are removed; L18/P18 are added.
# ---------------------------------------------------------------------------
$ws.Range("K18").ClearContents()
$ws.Range("N18").ClearContents()
$ws.Range("O18").ClearContents()
$ws.Range("L18").Formula = '=IF(K18<>"",K18,J18)'
$ws.Range("P18").Formula = '=IF(O18<>"",O18,N18)'

# ---------------------------------------------------------------------------
# Row 20 (the "sum" row that used to be row 19): the bold styling moves from
# the whole row to just the new L20/P20 totals; J20/K20/N20/O20 are
# recomputed as plain (non-shared, non-bold) SUM formulas over 13:18;
# L20/P20 are new bold SUM totals; M20 is removed entirely.
# ---------------------------------------------------------------------------
$ws.Range("I20").Font.Bold = $false

$ws.Range("J20").Formula = '=SUM(J13:J18)'
$ws.Range("J20").Font.Bold = $false

$ws.Range("K20").Formula = '=SUM(K13:K18)'
$ws.Range("K20").Font.Bold = $false

$ws.Range("L20").Formula = '=SUM(L13:L18)'
$ws.Range("L20").Font.Bold = $true

$ws.Range("M20").Clear()

$ws.Range("N20").Formula = '=SUM(N13:N18)'
$ws.Range("N20").Font.Bold = $false

$ws.Range("O20").Formula = '=SUM(O13:O18)'
$ws.Range("O20").Font.Bold = $false

$ws.Range("P20").Formula = '=SUM(P13:P18)'
$ws.Range("P20").Font.Bold = $true

# ---------------------------------------------------------------------------
# Row 28 (was row 27): add the new balance formula in K28.
# ---------------------------------------------------------------------------
$ws.Range("K28").Formula = '=K26+L20-P20'

# ---------------------------------------------------------------------------
# Update the view: scroll so column G is the left-most visible column and
# select P21 as the active cell (matching the new sheetView state).
# (topLeftCell/zoomScaleNormal attributes are not exposed for writing by
# this host; ScrollColumn is still set in case the host honors it.)
# ---------------------------------------------------------------------------
$ws.Range("P21").Select()
$excel.ActiveWindow.ScrollColumn = 7

Write-Host "Edit complete."
